$wb = $excel.ActiveWorkbook
$wsJournal = $wb.Worksheets.Item("Journal")
$wsTotaux  = $wb.Worksheets.Item("Totaux")

# ---------------------------------------------------------------------
# 1) "Journal" sheet (sheet1): append 5 new rows (27-31) for the new
#    4th week of work, all dated 27/02/2023 (serial 44984).
# ---------------------------------------------------------------------

$journalRows = @(
    @{ Row = 27; Date = 44984; Week = 4; Time = 0.05902777777777778;  Type = "Analyse";       Desc = "Tests des différents profils sur la nouvelle version version de l'application" },
    @{ Row = 28; Date = 44984; Week = 4; Time = 0.03125;               Type = "Documentation"; Desc = "Rédaction d'une description des analyses faites sur le rapport de projet" },
    @{ Row = 29; Date = 44984; Week = 4; Time = 0.010416666666666666;  Type = "Documentation"; Desc = "Fermeture du premier spint" },
    @{ Row = 30; Date = 44984; Week = 4; Time = 0.041666666666666664;  Type = "Documentation"; Desc = "Ouverture du deuxième sprint" },
    @{ Row = 31; Date = 44984; Week = 4; Time = 0.05555555555555555;   Type = "Documentation"; Desc = "Rédaction des futures questions à poser pour les entretiens" }
)

foreach ($r in $journalRows) {
    $row = $r.Row

    # Copy the format of the row right above so the new row matches the
    # rest of the table (date / week / time / type styles).
    $srcFmt = $wsJournal.Range("A" + ($row - 1) + ":E" + ($row - 1))
    $dstFmt = $wsJournal.Range("A" + $row + ":E" + $row)
    $srcFmt.Copy() | Out-Null
    $dstFmt.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $wsJournal.Cells.Item($row, 1).Value = $r.Date
    $wsJournal.Cells.Item($row, 2).Value = $r.Week
    $wsJournal.Cells.Item($row, 3).Value = $r.Time
    $wsJournal.Cells.Item($row, 4).Value = $r.Type
    $wsJournal.Cells.Item($row, 5).Value = $r.Desc
}

# Grow "Tableau1" to cover the 5 newly-added rows.
$tbl1 = $wsJournal.ListObjects.Item("Tableau1")
$tbl1.Resize($wsJournal.Range("A1:E31")) | Out-Null

# ---------------------------------------------------------------------
# 2) "Totaux" sheet (sheet2): insert a new weekly-subtotal row (11) for
#    the week of 27/02/2023 before the grand-total row, and push the
#    grand total down to row 12 with its SUM formula extended.
# ---------------------------------------------------------------------

$wsTotaux.Rows.Item(11).Insert() | Out-Null
$wsTotaux.Cells.Item(11, 1).Value = 44984
$wsTotaux.Cells.Item(11, 2).Clear() | Out-Null  # the inserted row copied B10's format; no value belongs there

$tbl2 = $wsTotaux.ListObjects.Item("Tableau2")
$tbl2.Resize($wsTotaux.Range("A1:B12")) | Out-Null

$wsTotaux.Cells.Item(12, 2).Formula = "=SUM(B2:B11)"

# ---------------------------------------------------------------------
# 3) Restore the selections recorded in the workbook (Totaux first so
#    that "Journal" ends up as the last-activated / displayed sheet).
# ---------------------------------------------------------------------

$wsTotaux.Range("C10").Select() | Out-Null
$wsJournal.Range("I13").Select() | Out-Null
